$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: update price (D) and 1h volume/change (E) columns,
# plus the OKB/ONDO rank swap (rows 44-45, columns B/C).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.300.50"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.353.78"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.46"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.37"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.18%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.349.31"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.470"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.47"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.93%  "
$ws.Range("E11").Value = "  -3.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.389"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.926.61"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.73%  "
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("E15").Value = "  -4.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.85"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.349.57"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.368.07"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.94"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.86"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.30"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "375.31"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.551"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.506.82"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.85"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("E27").Value = "  -2.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.79"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +7.69%  "
$ws.Range("E29").Value = "  -4.43%  "
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("E31").Value = "  +2.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.15"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.15"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.73%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.54"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("E36").Value = "  -6.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.80"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.53"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.76"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("E40").Value = "  -4.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("E42").Value = "  -1.23%  "
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.34"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.20"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.37"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.01"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.84"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.87"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.355.72"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("E51").Value = "  -2.65%  "
